$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.628.50'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '1.770.35'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''224.08'
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('D6').Value = '''0.544'
$ws.Range('E6').Value = '  -0.93%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '''31.74'
$ws.Range('E8').Value = '  +0.94%  '
$ws.Range('E9').Value = '  +1.42%  '
$ws.Range('E10').Value = '  -4.25%  '
$ws.Range('D11').Value = '''0.0934'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('D12').Value = '2.024.49'
$ws.Range('E12').Value = '  -0.91%  '
$ws.Range('D13').Value = '''11.01'
$ws.Range('E13').Value = '  +4.04%  '
$ws.Range('D14').Value = '1.773.50'
$ws.Range('E14').Value = '  -0.93%  '
$ws.Range('D15').Value = '33.667.95'
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('D16').Value = '''0.607'
$ws.Range('E16').Value = '  -3.22%  '
$ws.Range('E17').Value = '  -2.59%  '
$ws.Range('D18').Value = '''66.40'
$ws.Range('D19').Value = '0.0₃0774'
$ws.Range('E19').Value = '  -1.13%  '
$ws.Range('D20').Value = '''237.79'
$ws.Range('E20').Value = '  -2.88%  '
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('E22').Value = '  -1.93%  '
$ws.Range('D23').Value = '''3.99'
$ws.Range('E23').Value = '  -2.12%  '
$ws.Range('E24').Value = '  -2.64%  '
$ws.Range('E25').Value = '  +0.79%  '
$ws.Range('E26').Value = '  -1.82%  '
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('E28').Value = '  -0.26%  '
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('E30').Value = '  +1.33%  '
$ws.Range('D31').Value = '''0.0510'
$ws.Range('E31').Value = '  -1.42%  '
$ws.Range('D32').Value = '''3.59'
$ws.Range('E32').Value = '  -2.70%  '
$ws.Range('E33').Value = '  -0.31%  '
$ws.Range('E34').Value = '  -1.71%  '
$ws.Range('D35').Value = '1.379.43'
$ws.Range('E35').Value = '  -2.03%  '
$ws.Range('D36').Value = '''0.644'
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('D37').Value = '''1.03'
$ws.Range('E37').Value = '  -2.36%  '
$ws.Range('E38').Value = '  -1.39%  '
$ws.Range('D39').Value = '''2.23'
$ws.Range('E39').Value = '  +5.72%  '
$ws.Range('E40').Value = '  +0.81%  '
$ws.Range('D41').Value = '''77.81'
$ws.Range('E41').Value = '  -2.24%  '
$ws.Range('E42').Value = '  -3.79%  '
$ws.Range('E43').Value = '  -2.40%  '
$ws.Range('D44').Value = '''13.43'
$ws.Range('E44').Value = '  +13.54%  '
$ws.Range('E45').Value = '  +3.88%  '
$ws.Range('E46').Value = '  +13.84%  '
$ws.Range('D47').Value = '''0.0499'
$ws.Range('E47').Value = '  +0.98%  '
$ws.Range('E48').Value = '  +1.54%  '
$ws.Range('D49').Value = '''5.82'
$ws.Range('E49').Value = '  -2.11%  '
$ws.Range('D50').Value = '1.925.32'
$ws.Range('E50').Value = '  -0.53%  '
